{"js": "// Update the date paragraph and the 25 division-problem table cells to the\n// new values from the diff. Replacements are applied by position (not by\n// searching for old text), since some old/new values repeat or collide.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// First paragraph holds the date line.\nparagraphs.items[0].insertText(\"2024-04-30 Tuesday\", Word.InsertLocation.replace);\n\n// The table holds the practice problems. Content lives in every 4th row\n// (rows 0, 4, 8, 12, 16); the rows in between are empty spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\n// New values, in row-major reading order across the five populated rows.\nconst newValues = [\n  [\"23\u00f74=5, 3\", \"60\u00f72=30, 0\", \"85\u00f73=28, 1\", \"77\u00f77=11, 0\", \"34\u00f79=3, 7\"],\n  [\"79\u00f72=39, 1\", \"35\u00f75=7, 0\", \"85\u00f78=10, 5\", \"23\u00f74=5, 3\", \"69\u00f77=9, 6\"],\n  [\"85\u00f73=28, 1\", \"68\u00f74=17, 0\", \"41\u00f73=13, 2\", \"62\u00f74=15, 2\", \"33\u00f74=8, 1\"],\n  [\"86\u00f79=9, 5\", \"27\u00f76=4, 3\", \"27\u00f78=3, 3\", \"61\u00f75=12, 1\", \"32\u00f79=3, 5\"],\n  [\"48\u00f76=8, 0\", \"22\u00f75=4, 2\", \"24\u00f79=2, 6\", \"73\u00f74=18, 1\", \"94\u00f75=18, 4\"],\n];\n\nconst contentRowIndices = [0, 4, 8, 12, 16];\n\nconst cellParagraphs = [];\nfor (let r = 0; r < contentRowIndices.length; r++) {\n  const rowIndex = contentRowIndices[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    cell.body.paragraphs.load(\"items\");\n    cellParagraphs.push({ paragraphsProxy: cell.body.paragraphs, text: newValues[r][c] });\n  }\n}\nawait context.sync();\n\nfor (const entry of cellParagraphs) {\n  entry.paragraphsProxy.items[0].insertText(entry.text, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the date paragraph and the 25 division-problem table cells to the\n# new values from the diff. Replacements are applied by position (row/column\n# index and paragraph index), not by searching for old text, since some\n# old/new values repeat or collide with each other.\n\n$d = $word.ActiveDocument\n\n# First paragraph holds the date line.\n$d.Paragraphs.Item(1).Range.Text = \"2024-04-30 Tuesday\"\n\n# The table holds the practice problems. Content lives in every 4th row\n# (1-based rows 1, 5, 9, 13, 17); the rows in between are empty spacer rows.\n$t = $d.Tables.Item(1)\n\n$contentRows = @(1, 5, 9, 13, 17)\n\n$newValues = @(\n  @(\"23\u00f74=5, 3\", \"60\u00f72=30, 0\", \"85\u00f73=28, 1\", \"77\u00f77=11, 0\", \"34\u00f79=3, 7\"),\n  @(\"79\u00f72=39, 1\", \"35\u00f75=7, 0\", \"85\u00f78=10, 5\", \"23\u00f74=5, 3\", \"69\u00f77=9, 6\"),\n  @(\"85\u00f73=28, 1\", \"68\u00f74=17, 0\", \"41\u00f73=13, 2\", \"62\u00f74=15, 2\", \"33\u00f74=8, 1\"),\n  @(\"86\u00f79=9, 5\", \"27\u00f76=4, 3\", \"27\u00f78=3, 3\", \"61\u00f75=12, 1\", \"32\u00f79=3, 5\"),\n  @(\"48\u00f76=8, 0\", \"22\u00f75=4, 2\", \"24\u00f79=2, 6\", \"73\u00f74=18, 1\", \"94\u00f75=18, 4\")\n)\n\nfor ($r = 0; $r -lt $contentRows.Length; $r++) {\n  $row = $contentRows[$r]\n  for ($c = 1; $c -le 5; $c++) {\n    $t.Cell($row, $c).Range.Text = $newValues[$r][$c - 1]\n  }\n}\n"}
